$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Account number column (C) holds a 24-digit numeric-looking value that must
# stay text (it already does for C2; force the same storage for the other
# rows that now share that account number).
$ws.Range("C3:C5").NumberFormat = "@"

# Row 2: JEMAA HORMI (C2 already holds this account number as text; leave it untouched)
$ws.Range("A2").Value = "JEMAA HORMI"
$ws.Range("B2").Value = "B219321"
$ws.Range("D2").Value = "KHOURIBGA"
$ws.Range("E2").Value = "CA"
$ws.Range("F2").Value = "Direction régionale"
$ws.Range("G2").Value = "001/RRR/AV1"
$ws.Range("H2").Value = "mensuelle"
$ws.Range("I2").Value = 7000
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 7000

# Row 3: MOHAMED BADRANE
$ws.Range("A3").Value = "MOHAMED BADRANE"
$ws.Range("B3").Value = "I83603"
$ws.Range("C3").Value = "225400000805987601012173"
$ws.Range("D3").Value = "KHOURIBGA"
$ws.Range("E3").Value = "CA"
$ws.Range("F3").Value = "Direction régionale"
$ws.Range("G3").Value = "001/RRR/AV1"
$ws.Range("H3").Value = "mensuelle"
$ws.Range("I3").Value = 14000
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 14000

# Row 4: JEMAA HORMI (second transfer)
$ws.Range("A4").Value = "JEMAA HORMI"
$ws.Range("B4").Value = "B219321"
$ws.Range("C4").Value = "225400000805987601012173"
$ws.Range("D4").Value = "KHOURIBGA"
$ws.Range("E4").Value = "CA"
$ws.Range("F4").Value = "Direction régionale"
$ws.Range("G4").Value = "001/RRR/AV1"
$ws.Range("H4").Value = "mensuelle"
$ws.Range("I4").Value = 1000
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 1000

# Row 5: MOHAMED BADRANE (second transfer) - new row
$ws.Range("A5").Value = "MOHAMED BADRANE"
$ws.Range("B5").Value = "I83603"
$ws.Range("C5").Value = "225400000805987601012173"
$ws.Range("D5").Value = "KHOURIBGA"
$ws.Range("E5").Value = "CA"
$ws.Range("F5").Value = "Direction régionale"
$ws.Range("G5").Value = "001/RRR/AV1"
$ws.Range("H5").Value = "mensuelle"
$ws.Range("I5").Value = 2000
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 2000

# Row 6: new totals row (blank/space labels, summed amounts)
$ws.Range("A6").Value = " "
$ws.Range("B6").Value = " "
$ws.Range("C6").Value = " "
$ws.Range("D6").Value = " "
$ws.Range("E6").Value = " "
$ws.Range("F6").Value = " "
$ws.Range("G6").Value = " "
$ws.Range("H6").Value = " "
$ws.Range("I6").Value = 24000
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 24000
